$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "PROVINCIA" column (column C) entirely, shifting subsequent
# columns left. This corresponds to the commit "Remove Provincia Certificado".
$ws.Columns("C").Delete()

# Leave selection on D4 as in the final saved state.
$ws.Range("D4").Select()
